# Auto-generated: apply latest crypto price/volume snapshot updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.692.03"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "2.526.05"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'309.31"
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("D6").Value = "'100.34"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").Value = "'0.567"
$ws.Range("E7").Value = "  -1.19%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").Value = "'35.49"
$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "  -1.04%  "

$ws.Range("E12").Value = "  -2.47%  "

$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "2.917.42"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'15.31"
$ws.Range("E15").Value = "  -4.11%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.535.91"
$ws.Range("E16").Value = "  -2.80%  "

$ws.Range("D17").Value = "'0.812"
$ws.Range("E17").Value = "  -4.23%  "

$ws.Range("D18").Value = "42.682.73"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").Value = "'12.22"
$ws.Range("E21").Value = "  -3.02%  "

$ws.Range("D22").Value = "'69.27"
$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("D23").Value = "'242.81"
$ws.Range("E23").Value = "  -3.04%  "

$ws.Range("E24").Value = "  -3.10%  "

$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'25.35"
$ws.Range("E27").Value = "  -6.86%  "

$ws.Range("E28").Value = "  -2.11%  "

$ws.Range("D29").Value = "'10.15"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("D30").Value = "'38.47"
$ws.Range("E30").Value = "  -3.54%  "

$ws.Range("D31").Value = "'160.04"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").Value = "'5.77"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("E33").Value = "  +9.00%  "

$ws.Range("D34").Value = "'2.69"
$ws.Range("E34").Value = "  +0.97%  "

$ws.Range("D36").Value = "'18.41"
$ws.Range("E36").Value = "  -1.30%  "

$ws.Range("D37").Value = "'3.11"
$ws.Range("E37").Value = "  -7.39%  "

$ws.Range("E38").Value = "  -7.55%  "

$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").Value = "'4.17"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").Value = "'22.27"
$ws.Range("E42").Value = "  -4.45%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0299"
$ws.Range("E44").Value = "  -1.24%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.28"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").Value = "2.002.51"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").Value = "'8.84"
$ws.Range("E47").Value = "  -2.23%  "

$ws.Range("D48").Value = "2.769.80"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("E49").Value = "  -4.09%  "

$ws.Range("D50").Value = "'79.16"
$ws.Range("E50").Value = "  -3.61%  "

$ws.Range("D51").Value = "'100.65"
$ws.Range("E51").Value = "  -1.75%  "
